$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.967.22"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.82%  '
$ws.Range('D3').Value = "'1.871.00"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.02%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').Value = "'318.38"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.45%  '
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').Value = "'0.4355"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.52%  '
$ws.Range('E8').Value = '  -1.91%  '
$ws.Range('D9').Value = "'0.07478"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.69%  '
$ws.Range('D10').Value = "'0.9376"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.28%  '
$ws.Range('D11').Value = "'21.26"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.72%  '
$ws.Range('D12').Value = "'1.928.28"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.33%  '
$ws.Range('D13').Value = "'6.748"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.34%  '
$ws.Range('D14').Value = "'5.439"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.56%  '
$ws.Range('D15').Value = "'0.06869"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('D16').Value = "'1.004"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').Value = "'81.48"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.49%  '
$ws.Range('D18').Value = "'0.000009054"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.14%  '
$ws.Range('D19').Value = "'1.002"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.32%  '
$ws.Range('D20').Value = "'15.80"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.73%  '
$ws.Range('D21').Value = "'27.947.16"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.03%  '
$ws.Range('D22').Value = "'5.127"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.15%  '
$ws.Range('D23').Value = "'11.05"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('D24').Value = "'2.131.64"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.37%  '
$ws.Range('D25').Value = "'2.039"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.93%  '
$ws.Range('D26').Value = "'153.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.18%  '
$ws.Range('D27').Value = "'18.56"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.97%  '
$ws.Range('D28').Value = "'5.588"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.47%  '
$ws.Range('D29').Value = "'113.55"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.89%  '
$ws.Range('D30').Value = "'1.702"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -8.02%  '
$ws.Range('D31').Value = "'0.09024"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.41%  '
$ws.Range('D32').Value = "'0.8109"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.27%  '
$ws.Range('D33').Value = "'4.809"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.45%  '
$ws.Range('D34').Value = "'1.181"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.51%  '
$ws.Range('D35').Value = "'2.974"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.22%  '
$ws.Range('D36').Value = "'1.003"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('D37').Value = "'0.05515"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.37%  '
$ws.Range('D38').Value = "'1.120"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.38%  '
$ws.Range('D39').Value = "'0.01982"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.44%  '
$ws.Range('D40').Value = "'2.976"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.03%  '
$ws.Range('E41').Value = '  -4.58%  '
$ws.Range('D42').Value = "'0.1699"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.28%  '
$ws.Range('D43').Value = "'6.981"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.59%  '
$ws.Range('D44').Value = "'8.788"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.25%  '
$ws.Range('D45').Value = "'0.06755"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.52%  '
$ws.Range('D46').Value = "'0.4891"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.33%  '
$ws.Range('D47').Value = "'10.56"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.31%  '
$ws.Range('D48').Value = "'106.85"
$ws.Range('D48').Style = 'Normal'
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').Value = "'1.002"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.22%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = "'1.673"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.01%  '
$ws.Range('D51').Value = "'1.900"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -13.76%  '
